# "fixed ingredient and category bugs in product form"
#
# Applies the tracked changes to the "Feuil1" sheet:
#  - "liste des catégories" block (rows 43-47): ticks the first four phase
#    checkboxes for the ingredient/category related rows that are now done,
#    and records a "2j" time estimate for the "liste des catégories" row.
#  - "liste des ingrédients" row (38): time estimate corrected from
#    "1 semaine" to "2 sem".
#  - Removes two now-empty spacer rows (64-65) below the notes table.
#  - Restores the view: frozen-pane scroll position and active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- shared-string bearing "temps" (time estimate) cells -------------------
# Order matters: new shared strings are appended in first-write order, and
# the target file expects "2j" (I43) before "2 sem" (I38) in the table.
$ws.Range("I43").Value = "2j"
$ws.Range("I38").Value = "2 sem"

# --- "liste des catégories" sub-tasks: mark phases as completed ------------
$ws.Range("C43:F43").Value = $true
$ws.Range("C44:F44").Value = $true
$ws.Range("C45:F45").Value = $true
$ws.Range("C47:F47").Value = $true

# --- remove the two blank spacer rows (64 & 65) without shifting the rows
# below them (the "fournisseurs" block must stay anchored at rows 67-69) ---
$ws.Rows("64:65").Delete()
$ws.Rows("65:66").Insert()

# --- restore the view state ------------------------------------------------
$ws.Activate()
$ws.Range("A33").Select()
$ws.Range("I39").Select()
